$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue "D2" "61.658.30"
Set-TextValue "E2" "  -2.16%  "
Set-TextValue "D3" "3.404.77"
Set-TextValue "E3" "  -1.91%  "
Set-TextValue "E4" "  -0.11%  "
Set-TextValue "D5" "403.64"
Set-TextValue "E5" "  -1.01%  "
Set-TextValue "D6" "131.65"
Set-TextValue "E6" "  +0.14%  "
Set-TextValue "E7" "  -2.19%  "
Set-TextValue "E8" "  +0.02%  "
Set-TextValue "D9" "0.683"
Set-TextValue "E9" "  -1.81%  "
Set-TextValue "E10" "  -4.54%  "
Set-TextValue "D11" "41.71"
Set-TextValue "E11" "  -4.28%  "
Set-TextValue "E12" "  -1.13%  "
Set-TextValue "E13" "  -5.79%  "
Set-TextValue "D14" "19.72"
Set-TextValue "E14" "  -2.01%  "
Set-TextValue "D15" "3.410.33"
Set-TextValue "E15" "  -2.53%  "
Set-TextValue "D16" "11.61"
Set-TextValue "E16" "  +6.51%  "
Set-TextValue "D17" "61.640.21"
Set-TextValue "E17" "  -2.25%  "
Set-TextValue "E18" "  -3.57%  "
Set-TextValue "E19" "  -1.75%  "
Set-TextValue "D21" "82.98"
Set-TextValue "E21" "  +0.37%  "
Set-TextValue "D22" "311.05"
Set-TextValue "E22" "  -0.35%  "
Set-TextValue "D23" "12.68"
Set-TextValue "E23" "  -3.50%  "
Set-TextValue "D24" "3.13"
Set-TextValue "E24" "  -0.82%  "
Set-TextValue "E25" "  +10.24%  "
Set-TextValue "B26" "RenderToken"
Set-TextValue "C26" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D26" "8.08"
Set-TextValue "E26" "  +6.41%  "
Set-TextValue "B27" "EthereumClassic"
Set-TextValue "C27" "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue "D27" "8.01"
Set-TextValue "E27" "  -2.07%  "
Set-TextValue "D28" "8.01"
Set-TextValue "E28" "  -2.07%  "
Set-TextValue "D29" "2.75"
Set-TextValue "E29" "  +5.62%  "
Set-TextValue "D30" "0.171"
Set-TextValue "E30" "  -3.73%  "
Set-TextValue "E31" "  -2.77%  "
Set-TextValue "D32" "42.88"
Set-TextValue "E32" "  -3.25%  "
Set-TextValue "D33" "0.999"
Set-TextValue "E33" "  +0.11%  "
Set-TextValue "E34" "  -4.43%  "
Set-TextValue "D35" "0.0481"
Set-TextValue "E35" "  -2.76%  "
Set-TextValue "D36" "51.36"
Set-TextValue "E36" "  -2.36%  "
Set-TextValue "D37" "0.998"
Set-TextValue "E37" "  -0.18%  "
Set-TextValue "E38" "  -5.82%  "
Set-TextValue "B39" "Stacks"
Set-TextValue "C39" "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue "D39" "2.95"
Set-TextValue "E39" "  -2.61%  "
Set-TextValue "B40" "TheGraph"
Set-TextValue "C40" "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
Set-TextValue "D40" "0.320"
Set-TextValue "E40" "  +11.07%  "
Set-TextValue "D41" "139.30"
Set-TextValue "E41" "  +1.56%  "
Set-TextValue "E42" "  -1.40%  "
Set-TextValue "E43" "  -1.60%  "
Set-TextValue "D44" "3.94"
Set-TextValue "E44" "  -0.85%  "
Set-TextValue "D45" "16.51"
Set-TextValue "E45" "  -5.97%  "
Set-TextValue "D46" "2.22"
Set-TextValue "E46" "  -0.95%  "
Set-TextValue "D47" "21.14"
Set-TextValue "E47" "  -4.91%  "
Set-TextValue "D48" "2.096.75"
Set-TextValue "E48" "  -4.06%  "
Set-TextValue "E49" "  -1.98%  "
Set-TextValue "D50" "1.76"
Set-TextValue "E50" "  +20.06%  "
Set-TextValue "D51" "1.92"
Set-TextValue "E51" "  +2.55%  "
